$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily log. Insert a new row at
# position 386 (shifting the existing rows 386-435 down to 387-436) and
# populate it with the new observation.
$ws.Rows.Item(386).Insert()

$ws.Range("A386").Value = 8
$ws.Range("B386").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C386").Value = 'Coquimbo'
$ws.Range("D386").Value = 44984
$ws.Range("E386").Value = 4
$ws.Range("F386").Value = 100112003
$ws.Range("G386").Value = 'Ajo'
$ws.Range("H386").Value = 'Chino'
$ws.Range("I386").Value = 'Primera'
$ws.Range("J386").Value = 440
$ws.Range("K386").Value = 17000
$ws.Range("L386").Value = 18000
$ws.Range("M386").Value = 17500
$ws.Range("N386").Value = '$/caja 10 kilos'
$ws.Range("O386").Value = 'China'
$ws.Range("P386").Value = 1750
$ws.Range("Q386").Value = 10
$ws.Range("R386").Value = 'Hortaliza'
